# Veh Buyer Discount Rate.xlsx -- bring InputData up to eps v2.0.0
#
# Real content changes (everything else in the recorded diff is Excel's own
# save-time bookkeeping -- fileVersion/calcId/xr GUIDs/theme font cache/
# customXml item ids -- which isn't reachable, or meaningful, through the
# object model):
#   1. "About" sheet, A11: fix -> "Vehcile buyer discount rates ..."
#   2. "VBDR" sheet,  B1 : "Discount Rate" -> "Discount Rate (dimensionless)"
#   3. Restore the saved selection on each sheet (About -> A1, VBDR -> B2)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsVBDR  = $wb.Worksheets.Item("VBDR")

$wsAbout.Range("A11").Value = "Vehcile buyer discount rates vary tremendously by study."
$wsVBDR.Range("B1").Value = "Discount Rate (dimensionless)"

# Match the saved view state: VBDR selects B2, About keeps A1 as the
# active cell, and the About tab stays the selected/active one.
$wsVBDR.Range("B2").Select()
$wsAbout.Range("A1").Select()
$wsAbout.Activate()
